$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.736.78"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.894.17"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.26"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4925"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2953"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06811"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "WrappedEther"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.894.01"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.40"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.98%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "92.45"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07263"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6858"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.92%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.67%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "30.710.43"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.64%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.135.03"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9992"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "194.83"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +40.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.097"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +8.41%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.03"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +4.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.55"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +15.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.929"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.396"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.348"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +3.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09028"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.034"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +2.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05193"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7492"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +5.53%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.699"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01884"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +7.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.674"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9378"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4453"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.85"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.829"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.729"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1343"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +7.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05860"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +3.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.728"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.13%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.02%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "33.61"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +4.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.409"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.89%  "
